$wb = $excel.ActiveWorkbook

$targetFile = "89674ea8-2bd7-4120-a91e-5ee1e7dfcba4.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/f7cbb412f5b06a370152b85a6fe41bc04ca19fd4/e2e/89674ea8-2bd7-4120-a91e-5ee1e7dfcba4.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/62954174486738efdd6d1ccd37878b7da84b647d/e2e/89674ea8-2bd7-4120-a91e-5ee1e7dfcba4.md."
$currentTargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/f7cbb412f5b06a370152b85a6fe41bc04ca19fd4/e2e/89674ea8-2bd7-4120-a91e-5ee1e7dfcba4.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

# Widen the Error Detail column (P) to fit the new long message
$ws.Columns.Item(16).ColumnWidth = 39.17

# Row 8 corresponds to 89674ea8-2bd7-4120-a91e-5ee1e7dfcba4.md - fill in the handback-report columns
$ws.Range("I8").Value = $targetFile
$ws.Range("J8").Value = $ws.Range("G8").Value()
$ws.Range("K8").Value = "2016-08-13 12:54:27"
$ws.Range("P8").Value = $errorDetail

# Give I8 the same "hyperlink" look as the other linked cells, and wire the hyperlink itself
$ws.Range("I8").Font.Underline = 2
$ws.Range("I8").Font.Color = 15570276
$ws.Hyperlinks.Add($ws.Range("I8"), $currentTargetUrl, "", "", $targetFile)

# ---- de-de sheet ----
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Columns.Item(16).ColumnWidth = 39.17

$ws2.Range("I8").Value = $targetFile
$ws2.Range("J8").Value = $ws2.Range("G8").Value()
$ws2.Range("K8").Value = "2016-08-13 12:54:38"
$ws2.Range("P8").Value = $errorDetail

$ws2.Range("I8").Font.Underline = 2
$ws2.Range("I8").Font.Color = 15570276
$ws2.Hyperlinks.Add($ws2.Range("I8"), $currentTargetUrl, "", "", $targetFile)
